# Update the "取得日時" (acquisition timestamp) column for rows 2-8 on the
# "ランサーズ" sheet to reflect the new append timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-12 07:06:17"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
